# "ajout tOn et tOff" - add tOn/tOff/rOn/rOff columns (duty-cycle helper
# columns) between the existing "consoOff" and "tension" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new columns (tOn, tOff, rOn, rOff) right before the
# existing "tension" column (was D, becomes H).
$ws.Columns("D:G").Insert()

$ws.Range("D1").Value = "tOn"
$ws.Range("E1").Value = "tOff"
$ws.Range("F1").Value = "rOn"
$ws.Range("G1").Value = "rOff"

# tOn (seconds the module is "on") - 2s for everything except the data
# logger row which stays on for 30s.
$ws.Range("D2:D7").Value = 2
$ws.Range("D8").Value = 30
$ws.Range("D9").Value = 2

# tOff (seconds the module is "off") - 60s for everything except the data
# logger row which only wakes up once every 24h.
$ws.Range("E2:E7").Value = 60
$ws.Range("E8").Formula = "=24*60*60"
$ws.Range("E9").Value = 60

# rOn / rOff: fraction of time spent on / off, derived from tOn & tOff.
$ws.Range("F2").Formula = "=D2/(E2+D2)"
$ws.Range("G2").Formula = "=E2/(E2+D2)"
$ws.Range("F3:F9").Formula = "=D3/(E3+D3)"
$ws.Range("G3:G9").Formula = "=E3/(E3+D3)"

# Match the existing "Satisfaisant" (green) look of the table, now with a
# higher-precision number format for all the numeric columns.
$ws.Range("B2:H9").NumberFormat = "0.00000000000"

# Resize the new/affected columns to fit their content.
$ws.Range("B:C").ColumnWidth = 12.65625
$ws.Range("D:D").ColumnWidth = 13.59375
$ws.Range("E:E").ColumnWidth = 16.875
$ws.Range("F:H").ColumnWidth = 12.65625

# Leave a blank formatted cell below the table (scratch/working cell) and
# select it, as in the authored workbook.
$ws.Range("G13").NumberFormat = "0.00E+00"
[void]$ws.Range("G13").Select()
